$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3
$ws.Range("AG2").Value = 126
$ws.Range("AI2").Value = 17
$ws.Range("BC2").Value = 101
$ws.Range("G3").Value = 2.55
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 2.55
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.75
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 11
$ws.Range("AT3").Value = 3.75
$ws.Range("AY3").Value = 13
$ws.Range("Q4").Value = 1.88
$ws.Range("R4").Value = 1.93
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 11
$ws.Range("AM4").Value = 26
$ws.Range("AS4").Value = 201
$ws.Range("AT4").Value = 2.75
$ws.Range("AU4").Value = 8
$ws.Range("BC4").Value = 151
$ws.Range("G5").Value = 2.2
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 3.4
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 19
$ws.Range("AK5").Value = 34
$ws.Range("AM5").Value = 34
$ws.Range("AN5").Value = 4.33
$ws.Range("AX5").Value = 5
$ws.Range("H6").Value = 4.5
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.7
$ws.Range("Y6").Value = 9
$ws.Range("G7").Value = 1.42
$ws.Range("I7").Value = 8.5
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("W7").Value = 5
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 9
$ws.Range("AB7").Value = 34
$ws.Range("AF7").Value = 101
$ws.Range("AH7").Value = 17
$ws.Range("AJ7").Value = 26
$ws.Range("AN7").Value = 3.1
$ws.Range("AX7").Value = 9
$ws.Range("Q8").Value = 2.1
$ws.Range("R8").Value = 1.7
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("G11").Value = 2.4
$ws.Range("I11").Value = 2.7
$ws.Range("L11").Value = 3.25
$ws.Range("Q11").Value = 1.73
$ws.Range("R11").Value = 2.08
$ws.Range("W11").Value = 10
$ws.Range("X11").Value = 13
$ws.Range("AA11").Value = 19
$ws.Range("AC11").Value = 12
$ws.Range("AJ11").Value = 10
$ws.Range("AZ11").Value = 21
$ws.Range("N12").Value = 13
$ws.Range("G13").Value = 1.36
$ws.Range("H13").Value = 4.5
$ws.Range("N13").Value = 12
$ws.Range("Y13").Value = 9
$ws.Range("AC13").Value = 12
$ws.Range("AD13").Value = 9
$ws.Range("AG13").Value = 351
$ws.Range("AJ13").Value = 23
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 2.2
$ws.Range("G15").Value = 2.55
$ws.Range("I15").Value = 2.35
$ws.Range("J15").Value = 3.25
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 15
$ws.Range("O15").Value = 1.18
$ws.Range("P15").Value = 4.5
$ws.Range("Q15").Value = 1.65
$ws.Range("R15").Value = 2.2
$ws.Range("Z15").Value = 29
$ws.Range("AA15").Value = 21
$ws.Range("AB15").Value = 26
$ws.Range("AI15").Value = 13
$ws.Range("AK15").Value = 23
$ws.Range("AN15").Value = 5
$ws.Range("AO15").Value = 15
$ws.Range("G16").Value = 1.83
$ws.Range("I16").Value = 4.33
$ws.Range("J16").Value = 2.6
$ws.Range("L16").Value = 5
$ws.Range("W16").Value = 6
$ws.Range("X16").Value = 8
$ws.Range("AC16").Value = 7.5
$ws.Range("AF16").Value = 67
$ws.Range("AH16").Value = 11
$ws.Range("AK16").Value = 51
$ws.Range("AO16").Value = 10
$ws.Range("AQ16").Value = 34
$ws.Range("AY16").Value = 26
$ws.Range("G18").Value = 2.9
$ws.Range("I18").Value = 2.5
$ws.Range("K18").Value = 2.05
$ws.Range("L18").Value = 3.2
$ws.Range("N18").Value = 10
$ws.Range("O18").Value = 1.33
$ws.Range("P18").Value = 3.25
$ws.Range("Q18").Value = 2.08
$ws.Range("R18").Value = 1.73
$ws.Range("X18").Value = 13
$ws.Range("AC18").Value = 9
$ws.Range("AE18").Value = 15
$ws.Range("AJ18").Value = 10
$ws.Range("AX18").Value = 4.5
$ws.Range("AY18").Value = 15
$ws.Range("BA18").Value = 51
$ws.Range("BC18").Value = 201
$ws.Range("G20").Value = 2.85
$ws.Range("I20").Value = 2.37
$ws.Range("L20").Value = 2.95
$ws.Range("N20").Value = 6.9
$ws.Range("O20").Value = 1.34
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 1.72
$ws.Range("S20").Value = 1.39
$ws.Range("T20").Value = 2.77
$ws.Range("U20").Value = 1.78
$ws.Range("V20").Value = 1.93
$ws.Range("W20").Value = 8.75
$ws.Range("X20").Value = 14.5
$ws.Range("Z20").Value = 35
$ws.Range("AC20").Value = 6.9
$ws.Range("AD20").Value = 6.2
$ws.Range("AH20").Value = 7.6
$ws.Range("AI20").Value = 11.25
$ws.Range("AK20").Value = 24
$ws.Range("AL20").Value = 20
$ws.Range("AM20").Value = 30
$ws.Range("AT20").Value = 2.77
$ws.Range("AZ20").Value = 20
$ws.Range("BB20").Value = 80
$ws.Range("Q21").Value = 1.7
$ws.Range("R21").Value = 2.1
$ws.Range("G22").Value = 1.83
$ws.Range("I22").Value = 3.75
$ws.Range("Q22").Value = 1.9
$ws.Range("R22").Value = 1.9
$ws.Range("U22").Value = 1.83
$ws.Range("V22").Value = 1.83
$ws.Range("AE22").Value = 17
$ws.Range("AM22").Value = 41
$ws.Range("G23").Value = 2.9
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 2.25
$ws.Range("U23").Value = 1.73
$ws.Range("V23").Value = 2
$ws.Range("AH23").Value = 8.5
$ws.Range("AI23").Value = 11
$ws.Range("AJ23").Value = 9.5
$ws.Range("AO23").Value = 17
$ws.Range("AY23").Value = 12
$ws.Range("G25").Value = 1.26
$ws.Range("I25").Value = 9.5
$ws.Range("K25").Value = 2.57
$ws.Range("L25").Value = 7.8
$ws.Range("N25").Value = 14.7
$ws.Range("O25").Value = 1.17
$ws.Range("P25").Value = 5.05
$ws.Range("Q25").Value = 1.5
$ws.Range("R25").Value = 2.25
$ws.Range("S25").Value = 1.27
$ws.Range("T25").Value = 3.55
$ws.Range("X25").Value = 6.4
$ws.Range("Y25").Value = 9
$ws.Range("Z25").Value = 7.6
$ws.Range("AA25").Value = 10.75
$ws.Range("AD25").Value = 11
$ws.Range("AH25").Value = 27
$ws.Range("AI25").Value = 75
$ws.Range("AJ25").Value = 30
$ws.Range("AK25").Value = 300
$ws.Range("AN25").Value = 3.1
$ws.Range("AO25").Value = 5.3
$ws.Range("AX25").Value = 10
$ws.Range("AY25").Value = 55
$ws.Range("AZ25").Value = 45
